# Updated cryptos list on Tue Apr 18 22:23:04 UTC 2023 with GitHub Actions
#
# Refreshes the Price (column D) and Volume(1h) (column E) figures for the
# crypto list, and fixes the ordering of three rows (13-15) whose coins
# were re-ranked: Chainlink now leads, then Polkadot, then WrappedEther
# (each carrying its own refreshed link/price/volume values).
#
# Numeric-looking price strings (e.g. "343.03") are written with a leading
# apostrophe so Excel keeps them as text instead of silently converting them
# to numbers, matching the original inline-string cell contents.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.406.49"

$ws.Range("D3").Value = "2.096.49"
$ws.Range("E3").Value = "  +0.10%  "

$ws.Range("E4").Value = "  -0.77%  "

$ws.Range("D5").Value = "'343.03"
$ws.Range("E5").Value = "  -0.07%  "

$ws.Range("E6").Value = "  -0.63%  "

$ws.Range("D7").Value = "'0.5249"
$ws.Range("E7").Value = "  +1.59%  "

$ws.Range("D8").Value = "'0.4429"
$ws.Range("E8").Value = "  +1.09%  "

$ws.Range("D9").Value = "'54.44"
$ws.Range("E9").Value = "  +3.29%  "

$ws.Range("D10").Value = "'0.09371"
$ws.Range("E10").Value = "  +1.03%  "

$ws.Range("D11").Value = "'1.170"

$ws.Range("D12").Value = "'24.80"
$ws.Range("E12").Value = "  -0.32%  "

# Row 13: WrappedEther -> Chainlink
$ws.Range("B13").Value = "Chainlink"
$ws.Range("C13").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D13").Value = "'8.606"
$ws.Range("E13").Value = "  +4.04%  "

# Row 14: Chainlink -> Polkadot
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").Value = "'6.926"
$ws.Range("E14").Value = "  +2.57%  "

# Row 15: Polkadot -> WrappedEther
$ws.Range("B15").Value = "WrappedEther"
$ws.Range("C15").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D15").Value = "2.032.73"
$ws.Range("E15").Value = "  -3.39%  "

$ws.Range("D16").Value = "'101.56"
$ws.Range("E16").Value = "  +2.08%  "

$ws.Range("D17").Value = "'0.00001160"
$ws.Range("E17").Value = "  +0.58%  "

$ws.Range("D18").Value = "'1.002"
$ws.Range("E18").Value = "  -0.73%  "

$ws.Range("E19").Value = "  +2.10%  "

$ws.Range("D20").Value = "'0.06686"
$ws.Range("E20").Value = "  +0.71%  "

$ws.Range("D21").Value = "'6.332"
$ws.Range("E21").Value = "  +2.27%  "

$ws.Range("D22").Value = "'1.002"
$ws.Range("E22").Value = "  -0.60%  "

$ws.Range("D23").Value = "30.421.49"
$ws.Range("E23").Value = "  +2.27%  "

$ws.Range("D24").Value = "'12.54"
$ws.Range("E24").Value = "  +0.38%  "

$ws.Range("D25").Value = "'2.310"
$ws.Range("E25").Value = "  -0.48%  "

$ws.Range("E26").Value = "  -0.24%  "

$ws.Range("D27").Value = "'163.09"
$ws.Range("E27").Value = "  +1.14%  "

$ws.Range("D28").Value = "'6.805"
$ws.Range("E28").Value = "  +8.68%  "

$ws.Range("D29").Value = "'2.511"
$ws.Range("E29").Value = "  +0.01%  "

$ws.Range("D30").Value = "'133.53"
$ws.Range("E30").Value = "  +0.43%  "

$ws.Range("E31").Value = "  +0.48%  "

$ws.Range("D32").Value = "'1.665"
$ws.Range("E32").Value = "  +0.83%  "

$ws.Range("E33").Value = "  -0.04%  "

$ws.Range("D34").Value = "'6.275"
$ws.Range("E34").Value = "  +1.94%  "

$ws.Range("D35").Value = "'3.876"
$ws.Range("E35").Value = "  -1.61%  "

$ws.Range("E36").Value = "  -0.34%  "

$ws.Range("D37").Value = "'0.02637"
$ws.Range("E37").Value = "  +2.44%  "

$ws.Range("D38").Value = "'0.06821"
$ws.Range("E38").Value = "  +1.63%  "

$ws.Range("D39").Value = "'0.7006"
$ws.Range("E39").Value = "  +1.65%  "

$ws.Range("D40").Value = "'12.61"
$ws.Range("E40").Value = "  +1.22%  "

$ws.Range("D41").Value = "'1.341"
$ws.Range("E41").Value = "  +1.80%  "

$ws.Range("D42").Value = "'0.2219"
$ws.Range("E42").Value = "  -0.19%  "

$ws.Range("D43").Value = "'0.6845"
$ws.Range("E43").Value = "  +1.09%  "

$ws.Range("D44").Value = "'14.38"
$ws.Range("E44").Value = "  +0.66%  "

$ws.Range("D45").Value = "'2.348"
$ws.Range("E45").Value = "  +1.27%  "

$ws.Range("D46").Value = "'1.002"
$ws.Range("E46").Value = "  -0.56%  "

$ws.Range("D47").Value = "'1.385"
$ws.Range("E47").Value = "  +19.32%  "

$ws.Range("D48").Value = "'3.637"
$ws.Range("E48").Value = "  +0.60%  "

$ws.Range("D49").Value = "'0.00000000352"
$ws.Range("E49").Value = "  -2.75%  "

$ws.Range("D50").Value = "'1.232"
$ws.Range("E50").Value = "  +9.93%  "

$ws.Range("D51").Value = "'1.219"
$ws.Range("E51").Value = "  -0.01%  "
